$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Range("D2:D15").Value = "Y"
$ws.Activate()
$ws.Range("D2:D15").Select()
